$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (H) to hold
# the new "property_category" field, shifting date/legislator_name/
# legislator_id one column to the right.
$ws.Columns("H:H").Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"
